$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 76.55956557647782
$ws.Range("B3").Value = 82.29120972311028
$ws.Range("B4").Value = 86.65350025614113
$ws.Range("I5").Value = 96.33318055674334
$ws.Range("I6").Value = 96.07810130999621
$ws.Range("I7").Value = 96.38607941810061
$ws.Range("C8").Value = 96.65434250808033
$ws.Range("C9").Value = 94.53623463623801
$ws.Range("C10").Value = 95.96361525735652
$ws.Range("D11").Value = 99.01316026831296
$ws.Range("D12").Value = 99.07113843333264
$ws.Range("D13").Value = 99.00021132898752
$ws.Range("E14").Value = 98.71186137178094
$ws.Range("E15").Value = 98.597013617126
$ws.Range("E16").Value = 98.70885762399384
$ws.Range("F17").Value = 98.45377569191345
$ws.Range("F18").Value = 98.59182803612681
$ws.Range("F19").Value = 98.44737085252088
$ws.Range("G20").Value = 97.97098866908939
$ws.Range("G21").Value = 97.93207065081768
$ws.Range("G22").Value = 97.81467016817102
$ws.Range("H23").Value = 96.29443142206942
$ws.Range("H24").Value = 96.89925224075718
$ws.Range("H25").Value = 97.26449572742595
$ws.Range("B26").Value = 82.72088649806335
$ws.Range("B27").Value = 87.57603012161735
$ws.Range("I28").Value = 95.95919307509584
$ws.Range("I29").Value = 96.13944521742309
$ws.Range("C30").Value = 96.40546529413945
$ws.Range("C31").Value = 95.51310931412293
$ws.Range("D32").Value = 98.98887047674418
$ws.Range("D33").Value = 99.06537830278359
$ws.Range("E34").Value = 98.53469641454124
$ws.Range("E35").Value = 98.38715958375845
$ws.Range("F36").Value = 98.56645294065849
$ws.Range("F37").Value = 98.5004755496111
$ws.Range("G38").Value = 97.90113280940443
$ws.Range("G39").Value = 97.83504888958186
$ws.Range("H40").Value = 98.31096540891743
$ws.Range("H41").Value = 97.62570206106432
$ws.Range("B42").Value = 85.5331246172685
$ws.Range("B43").Value = 90.53878617918024
$ws.Range("I44").Value = 96.3237088487616
$ws.Range("I45").Value = 95.79920077277052
$ws.Range("C46").Value = 96.96186648268487
$ws.Range("C47").Value = 97.06071081289554
$ws.Range("D48").Value = 99.02056470920523
$ws.Range("D49").Value = 98.8202905451746
$ws.Range("E50").Value = 98.67905778784485
$ws.Range("E51").Value = 98.65735412219536
$ws.Range("F52").Value = 98.46302556467661
$ws.Range("F53").Value = 98.27826906159287
$ws.Range("G54").Value = 97.71922575433993
$ws.Range("G55").Value = 97.81623614999057
$ws.Range("H56").Value = 97.09262591308647
$ws.Range("H57").Value = 97.38022716516274
